$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3252.3333
$ws.Range("J17").Value = 3252.3333
$ws.Range("L17").Value = 9756.999899999999
$ws.Range("N17").Value = -10092.9999

$ws.Range("H32").Value = 8864.333000000001
$ws.Range("I32").Value = 7230.5
$ws.Range("J32").Value = 9681.25
$ws.Range("K32").Value = 7230.5
$ws.Range("L32").Value = 9681.25
$ws.Range("M32").Value = -6904.5
$ws.Range("N32").Value = -10333.25

$ws.Range("H62").Value = 5500
$ws.Range("I62").Value = 5500
$ws.Range("K62").Value = 5500
$ws.Range("M62").Value = -4876

$ws.Range("H65").Value = 5500
$ws.Range("I65").Value = 5500
$ws.Range("K65").Value = 27500
$ws.Range("M65").Value = -24380

$ws.Range("H70").Value = 2066.5833
$ws.Range("I70").Value = 1900
$ws.Range("K70").Value = 5700
$ws.Range("M70").Value = -5430

$ws.Range("H73").Value = 2066.5833
$ws.Range("I73").Value = 1900
$ws.Range("K73").Value = 5700
$ws.Range("M73").Value = -4764

$ws.Range("H96").Value = 515.4545000000001
$ws.Range("I96").Value = 451.625
$ws.Range("J96").Value = 685.6667
$ws.Range("K96").Value = 1354.875
$ws.Range("L96").Value = 2057.0001
$ws.Range("M96").Value = 18.125
$ws.Range("N96").Value = -4803.0001

$ws.Range("H103").Value = 522.5714
$ws.Range("I103").Value = 86.666664
$ws.Range("J103").Value = 849.5
$ws.Range("K103").Value = 259.999992
$ws.Range("L103").Value = 2548.5
$ws.Range("M103").Value = 326.000008
$ws.Range("N103").Value = -3720.5

$ws.Range("H107").Value = 1092.625
$ws.Range("I107").Value = 679.61536
$ws.Range("J107").Value = 2882.3333
$ws.Range("K107").Value = 679.61536
$ws.Range("L107").Value = 2882.3333
$ws.Range("M107").Value = 1240.38464
$ws.Range("N107").Value = -6722.3333

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H132").Value = 4526.778
$ws.Range("I132").Value = 837.8
$ws.Range("K132").Value = 2513.4
$ws.Range("M132").Value = 16.60000000000036

$ws.Range("H135").Value = 1761.6
$ws.Range("I135").Value = 1109.7
$ws.Range("J135").Value = 3065.4
$ws.Range("K135").Value = 9987.300000000001
$ws.Range("L135").Value = 27588.6
$ws.Range("M135").Value = -7452.300000000001
$ws.Range("N135").Value = -32658.6

$ws.Range("H137").Value = 4047
$ws.Range("I137").Value = 4000
$ws.Range("J137").Value = 4049.9375
$ws.Range("K137").Value = 12000
$ws.Range("L137").Value = 12149.8125
$ws.Range("M137").Value = -9450
$ws.Range("N137").Value = -17249.8125

$ws.Range("H138").Value = 6727.3
$ws.Range("J138").Value = 6919.222
$ws.Range("L138").Value = 20757.666
$ws.Range("N138").Value = -31037.666

$ws.Range("H141").Value = 891
$ws.Range("I141").Value = 891
$ws.Range("K141").Value = 2673
$ws.Range("M141").Value = 2507

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2681.861
$ws.Range("I45").Value = 2504.4707
$ws.Range("K45").Value = 2504.4707
$ws.Range("M45").Value = -2127.4707

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 17799.75
$ws.Range("I54").Value = 12066.333
$ws.Range("J54").Value = 35000
$ws.Range("K54").Value = 12066.333
$ws.Range("L54").Value = 35000
$ws.Range("M54").Value = -11582.333
$ws.Range("N54").Value = -35968

$ws.Range("H86").Value = 6974.875
$ws.Range("J86").Value = 6974.875
$ws.Range("L86").Value = 6974.875
$ws.Range("N86").Value = -9220.875

$ws.Range("H89").Value = 6974.875
$ws.Range("J89").Value = 6974.875
$ws.Range("L89").Value = 34874.375
$ws.Range("N89").Value = -46106.375

$ws.Range("H94").Value = 1721.6
$ws.Range("I94").Value = 836
$ws.Range("K94").Value = 836
$ws.Range("M94").Value = -385

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32831.832
$ws.Range("I31").Value = 21373.625
$ws.Range("K31").Value = 21373.625
$ws.Range("M31").Value = -21078.625

$ws.Range("H34").Value = 32831.832
$ws.Range("I34").Value = 21373.625
$ws.Range("K34").Value = 21373.625
$ws.Range("M34").Value = -21171.625

$ws.Range("H41").Value = 19000
$ws.Range("I41").Value = 15000
$ws.Range("K41").Value = 15000
$ws.Range("M41").Value = -14572

$ws.Range("H60").Value = 18139
$ws.Range("J60").Value = 25000
$ws.Range("L60").Value = 25000
$ws.Range("N60").Value = -26022

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2351.5
$ws.Range("J34").Value = 3134.6667
$ws.Range("L34").Value = 9404.000100000001
$ws.Range("N34").Value = -9572.000100000001

$ws.Range("H55").Value = 5333.3335
$ws.Range("J55").Value = 5333.3335
$ws.Range("L55").Value = 16000.0005
$ws.Range("N55").Value = -16354.0005

$ws.Range("H100").Value = 1666.6666
$ws.Range("I100").Value = 1666.6666
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 4999.9998
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -4188.9998
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 4500
$ws.Range("I56").Value = 4500
$ws.Range("K56").Value = 4500
$ws.Range("M56").Value = -3809

$ws.Range("H82").Value = 2109.5
$ws.Range("I82").Value = 2244.3333
$ws.Range("J82").Value = 896
$ws.Range("K82").Value = 2244.3333
$ws.Range("L82").Value = 896
$ws.Range("M82").Value = -1883.3333
$ws.Range("N82").Value = -1618

$ws.Range("H85").Value = 2109.5
$ws.Range("I85").Value = 2244.3333
$ws.Range("J85").Value = 896
$ws.Range("K85").Value = 2244.3333
$ws.Range("L85").Value = 896
$ws.Range("M85").Value = -996.3332999999998
$ws.Range("N85").Value = -3392

$ws.Range("H106").Value = 6249.25
$ws.Range("J106").Value = 6249.25
$ws.Range("L106").Value = 6249.25
$ws.Range("N106").Value = -8773.25

$ws.Range("H122").Value = 5889.5835
$ws.Range("I122").Value = 6021.875
$ws.Range("J122").Value = 5625
$ws.Range("K122").Value = 18065.625
$ws.Range("L122").Value = 16875
$ws.Range("M122").Value = -15615.625
$ws.Range("N122").Value = -21775

$ws.Range("H132").Value = 45999.6
$ws.Range("I132").Value = 44999.5
$ws.Range("K132").Value = 134998.5
$ws.Range("M132").Value = -132468.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 30003
$ws.Range("I3").Value = 30003
$ws.Range("K3").Value = 30003
$ws.Range("M3").Value = -29889

$ws.Range("H52").Value = 6678680.5
$ws.Range("I52").Value = 6678680.5
$ws.Range("K52").Value = 6678680.5
$ws.Range("M52").Value = -6678454.5

$ws.Range("H81").Value = 433.6
$ws.Range("I81").Value = 429.33334
$ws.Range("K81").Value = 858.66668
$ws.Range("M81").Value = 202.33332

$ws.Range("H84").Value = 433.6
$ws.Range("I84").Value = 429.33334
$ws.Range("K84").Value = 4293.3334
$ws.Range("M84").Value = 1010.6666

$ws.Range("H96").Value = 825
$ws.Range("I96").Value = 800
$ws.Range("K96").Value = 800
$ws.Range("M96").Value = 573

$ws.Range("H113").Value = 675.6
$ws.Range("I113").Value = 644.5833
$ws.Range("K113").Value = 1933.7499
$ws.Range("M113").Value = 236.2501

$ws.Range("H132").Value = 4679.2593
$ws.Range("I132").Value = 4106.409
$ws.Range("J132").Value = 7199.8
$ws.Range("K132").Value = 12319.227
$ws.Range("L132").Value = 21599.4
$ws.Range("M132").Value = -9789.226999999999
$ws.Range("N132").Value = -26659.4
